$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Last status check on" timestamp in F1
$ws.Range("F1").Value = "Last status check on: 09.02.2022 03:00"

# D9: change from text "+0.6" to numeric 0.6
$ws.Range("D9").Value = 0.6

# E9: change from text "2022-02-09 02:46:44" to numeric date serial, formatted like the other date cells
$ws.Range("E9").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$d = Get-Date -Year 2022 -Month 2 -Day 9 -Hour 2 -Minute 46 -Second 44
$ws.Range("E9").Value = $d
